$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "004279859 / ASSOCIACAO / 1352225.88" record (row 2) ---
$ws.Rows.Item(2).Delete()

# --- Add a new record for "004813088 / JULIANA / 591.63" ---
# It belongs right after "004563252 / FERNANDO / 758.03", which - once the
# row above has been removed - now sits on row 36, so the new record goes
# on row 37 (pushing "004374891 / RODRIGO / 554.85" and everything below it
# down by one row).
$ws.Rows.Item(37).Insert()

# Account numbers are stored as text (leading zeros must be preserved).
# Enter it as text via a leading apostrophe, then strip the resulting
# "number stored as text" formatting so the cell keeps the plain, unstyled
# look used by every other account-number cell in the column.
$ws.Range("A37").Value = "'004813088"
$ws.Range("A37").ClearFormats()

$ws.Range("B37").Value = "JULIANA"
$ws.Range("C37").Value = 591.63
